## Generate Report for Handoff
## Rolls the localization-status report forward to the new source file
## (51f6b4a9-35f4-41bf-a8b6-f2d1334d6556.md replacing
## d4914487-b880-4503-a38d-9885374630c7.md), refreshes the handoff
## timestamps, and clears the stale handback info for the two language
## sheets (no handback has happened yet for the new file revision).

$wb = $excel.ActiveWorkbook

$oldGuid = "d4914487-b880-4503-a38d-9885374630c7"
$newGuid = "51f6b4a9-35f4-41bf-a8b6-f2d1334d6556"

# -------------------------------------------------------------------
# Overview sheet
# -------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "$newGuid.md"

# B2 keeps its original external link target, only the displayed path
# text is refreshed to the new file name.
$ovUrl = ""
foreach ($h in $ov.Hyperlinks) {
    $ovUrl = $h.Address
}
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), $ovUrl, $null, $null, "e2e\$newGuid.md")

$ov.Range("G2").Value = "2016-11-08 23:34:00"

# -------------------------------------------------------------------
# zh-cn sheet
# -------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zhUrlSource = ""
foreach ($h in $zh.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $zhUrlSource = $h.Address
    }
}
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $zhUrlSource, $null, $null, "$newGuid.md")

$zh.Range("G2").Value = "$newGuid.4bcce97dc34b0be4ac545ad182ac1452b031aec1.zh-cn.xlf"
$zh.Range("H2").Value = "2016-11-08 23:33:47"

# Latest Target File / Latest Handback File: no handback yet for the
# new revision, so both go blank and lose the old hyperlink formatting.
$zh.Range("I2").Value = "'"
$zh.Range("I2").Style = "Normal"
$zh.Range("J2").Value = "'"
$zh.Range("J2").Style = "Normal"

$zh.Range("K2").Value = "0001-01-01 00:00:00"

# Has metadata -> True (entered as quoted text so it isn't coerced to a
# native boolean, matching the original "True"/"False" text cells).
$zh.Range("O2").Value = "'True"
$zh.Range("O2").Style = "Normal"

$zh.Columns.Item(9).ColumnWidth = 17.75
$zh.Columns.Item(10).ColumnWidth = 20.75

# -------------------------------------------------------------------
# de-de sheet
# -------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$deUrlSource = ""
foreach ($h in $de.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $deUrlSource = $h.Address
    }
}
$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $deUrlSource, $null, $null, "$newGuid.md")

$de.Range("G2").Value = "$newGuid.4bcce97dc34b0be4ac545ad182ac1452b031aec1.de-de.xlf"
$de.Range("H2").Value = "2016-11-08 23:34:00"

$de.Range("I2").Value = "'"
$de.Range("I2").Style = "Normal"
$de.Range("J2").Value = "'"
$de.Range("J2").Style = "Normal"

$de.Range("K2").Value = "0001-01-01 00:00:00"

$de.Range("O2").Value = "'True"
$de.Range("O2").Style = "Normal"

$de.Columns.Item(9).ColumnWidth = 17.75
$de.Columns.Item(10).ColumnWidth = 20.75
